$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D2:E51 to text format first so numeric-looking strings (e.g. "1.00")
# are preserved exactly as text, matching the original inlineStr cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.739.80"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "3.332.48"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "580.89"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "175.54"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").Value = "3.328.72"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "46.49"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "704.66"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "3.878.38"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "67.731.02"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "3.335.84"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "17.39"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").Value = "11.02"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  +4.38%  "
$ws.Range("D24").Value = "16.97"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "98.54"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").Value = "9.49"
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").Value = "33.14"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "8.55"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "7.09"
$ws.Range("E31").Value = "  +4.85%  "
$ws.Range("D32").Value = "572.48"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").Value = "10.99"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").Value = "57.53"
$ws.Range("E35").Value = "  +3.83%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "3.711.58"
$ws.Range("E37").Value = "  -5.09%  "
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("E39").Value = "  +5.01%  "
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").Value = "2.64"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "0.0₃0676"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").Value = "0.0407"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Value = "2.68"
$ws.Range("E47").Value = "  +5.10%  "
$ws.Range("D48").Value = "0.129"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D50").Value = "1.33"
$ws.Range("E50").Value = "  -4.84%  "
$ws.Range("D51").Value = "128.93"
$ws.Range("E51").Value = "  +0.27%  "
